# Regenerate s_vals data to filter save games.
# Updates the numeric columns B:G for rows 2-6 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    3 = @(0.6753301551942219, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 1.642425054193055)
    4 = @(0.01514828764759746, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 0, 1.35982162114495)
    5 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 1, 12.59312877619104)
    6 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
